$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.395.22"
$ws.Range("E2").Value = "  -3.86%  "
$ws.Range("D3").Value = "3.034.41"
$ws.Range("E3").Value = "  -3.26%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "538.20"
$ws.Range("E5").Value = "  -4.51%  "
$ws.Range("D6").Value = "132.38"
$ws.Range("E6").Value = "  -10.66%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.026.51"
$ws.Range("E8").Value = "  -3.27%  "
$ws.Range("D9").Value = "0.483"
$ws.Range("E9").Value = "  -2.75%  "
$ws.Range("D10").Value = "6.36"
$ws.Range("E10").Value = "  -8.19%  "
$ws.Range("D11").Value = "0.153"
$ws.Range("E11").Value = "  -2.65%  "
$ws.Range("D12").Value = "0.452"
$ws.Range("E12").Value = "  -1.94%  "
$ws.Range("D13").Value = "34.29"
$ws.Range("E13").Value = "  -4.26%  "
$ws.Range("D14").Value = "0.0000211"
$ws.Range("E14").Value = "  -4.44%  "
$ws.Range("D15").Value = "3.533.99"
$ws.Range("E15").Value = "  -2.92%  "
$ws.Range("D16").Value = "62.539.83"
$ws.Range("E16").Value = "  -3.67%  "
$ws.Range("E17").Value = "  -1.87%  "
$ws.Range("D18").Value = "3.054.38"
$ws.Range("E18").Value = "  -2.65%  "
$ws.Range("D19").Value = "6.51"
$ws.Range("E19").Value = "  -3.07%  "
$ws.Range("D20").Value = "475.18"
$ws.Range("E20").Value = "  -9.79%  "
$ws.Range("D21").Value = "13.19"
$ws.Range("E21").Value = "  -4.37%  "
$ws.Range("D22").Value = "0.686"
$ws.Range("E22").Value = "  -1.72%  "
$ws.Range("D23").Value = "6.94"
$ws.Range("E23").Value = "  -6.24%  "
$ws.Range("D24").Value = "76.62"
$ws.Range("E24").Value = "  -2.38%  "
$ws.Range("D25").Value = "11.99"
$ws.Range("E25").Value = "  -5.39%  "
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").Value = "2.67"
$ws.Range("E27").Value = "  -4.12%  "
$ws.Range("D28").Value = "8.11"
$ws.Range("E28").Value = "  -5.35%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").Value = "1.90"
$ws.Range("E30").Value = "  -10.29%  "
$ws.Range("D31").Value = "25.78"
$ws.Range("E31").Value = "  -1.22%  "
$ws.Range("D32").Value = "1.11"
$ws.Range("E32").Value = "  -3.86%  "
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").Value = "59.82"
$ws.Range("E33").Value = "  +13.31%  "
$ws.Range("B34").Value = "Stacks"
$ws.Range("C34").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D34").Value = "2.44"
$ws.Range("E34").Value = "  -7.38%  "
$ws.Range("D35").Value = "507.74"
$ws.Range("E35").Value = "  -8.67%  "
$ws.Range("D36").Value = "5.84"
$ws.Range("E36").Value = "  -2.47%  "
$ws.Range("D37").Value = "5.03"
$ws.Range("E37").Value = "  -6.21%  "
$ws.Range("D38").Value = "0.0391"
$ws.Range("E38").Value = "  -10.00%  "
$ws.Range("D39").Value = "3.047.53"
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").Value = "0.0775"
$ws.Range("E40").Value = "  -4.53%  "
$ws.Range("D41").Value = "0.115"
$ws.Range("E41").Value = "  -4.23%  "
$ws.Range("D42").Value = "7.92"
$ws.Range("E42").Value = "  -3.60%  "
$ws.Range("D43").Value = "2.56"
$ws.Range("E43").Value = "  -9.41%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "0.247"
$ws.Range("E45").Value = "  -3.08%  "
$ws.Range("D46").Value = "2.00"
$ws.Range("E46").Value = "  -6.98%  "
$ws.Range("D47").Value = "119.66"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").Value = "23.83"
$ws.Range("E48").Value = "  -4.66%  "
$ws.Range("D49").Value = "0.105"
$ws.Range("E49").Value = "  -2.68%  "
$ws.Range("D50").Value = "0.0₃0485"
$ws.Range("E50").Value = "  -6.98%  "
$ws.Range("D51").Value = "2.30"
$ws.Range("E51").Value = "  +58.91%  "
